$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "309.67", "0.9980").
# Force the whole column range to Text format first so Excel stores the new
# values as text (preserving exact formatting like trailing zeros) instead of
# silently converting them to numbers, then clear the format change back off
# so the cells keep their original (default) style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.239.82"
$ws.Range("D3").Value = "1.683.62"
$ws.Range("D5").Value = "309.67"
$ws.Range("D6").Value = "0.9981"
$ws.Range("D7").Value = "0.3749"
$ws.Range("D8").Value = "0.3458"
$ws.Range("D11").Value = "0.07332"
$ws.Range("D12").Value = "0.9998"
$ws.Range("D14").Value = "6.144"
$ws.Range("D15").Value = "6.802"
$ws.Range("D16").Value = "1.683.80"
$ws.Range("D17").Value = "0.00001113"
$ws.Range("D18").Value = "0.9980"
$ws.Range("D19").Value = "0.06735"
$ws.Range("D20").Value = "82.51"
$ws.Range("D22").Value = "6.128"
$ws.Range("D23").Value = "12.08"
$ws.Range("D24").Value = "24.213.36"
$ws.Range("D26").Value = "2.698"
$ws.Range("D27").Value = "3.365"
$ws.Range("D28").Value = "152.91"
$ws.Range("D30").Value = "1.865.47"
$ws.Range("D31").Value = "127.39"
$ws.Range("D32").Value = "6.495"
$ws.Range("D33").Value = "4.072"
$ws.Range("D34").Value = "0.9967"
$ws.Range("D35").Value = "1.791"
$ws.Range("D36").Value = "0.08524"
$ws.Range("D37").Value = "12.61"
$ws.Range("D38").Value = "0.06511"
$ws.Range("D39").Value = "5.419"
$ws.Range("D40").Value = "9.009"
$ws.Range("D41").Value = "0.02363"
$ws.Range("D42").Value = "1.284"
$ws.Range("D43").Value = "0.2153"
$ws.Range("D44").Value = "0.6263"
$ws.Range("D45").Value = "0.9971"
$ws.Range("D46").Value = "13.35"
$ws.Range("D47").Value = "3.815"
$ws.Range("D48").Value = "0.6000"
$ws.Range("D49").Value = "128.20"
$ws.Range("D50").Value = "2.047"
$ws.Range("D51").Value = "0.07191"

$ws.Range("D2:D51").ClearFormats()

# Column E (percentages) is stored as plain text already (leading/trailing
# spaces keep Excel from parsing it as a number), so a direct assignment is fine.
$ws.Range("E2").Value = "  +11.51%  "
$ws.Range("E3").Value = "  +7.00%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("E5").Value = "  +8.66%  "
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("E8").Value = "  +5.77%  "
$ws.Range("E9").Value = "  +16.27%  "
$ws.Range("E10").Value = "  +5.62%  "
$ws.Range("E11").Value = "  +3.97%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("E14").Value = "  +5.70%  "
$ws.Range("E15").Value = "  +5.01%  "
$ws.Range("E16").Value = "  +7.66%  "
$ws.Range("E17").Value = "  +4.01%  "
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("E19").Value = "  +9.04%  "
$ws.Range("E20").Value = "  +11.66%  "
$ws.Range("E21").Value = "  +3.86%  "
$ws.Range("E22").Value = "  +5.04%  "
$ws.Range("E23").Value = "  +4.29%  "
$ws.Range("E24").Value = "  +11.45%  "
$ws.Range("E25").Value = "  +3.93%  "
$ws.Range("E26").Value = "  +11.92%  "
$ws.Range("E27").Value = "  -8.76%  "
$ws.Range("E28").Value = "  +3.04%  "
$ws.Range("E29").Value = "  +8.58%  "
$ws.Range("E30").Value = "  +7.46%  "
$ws.Range("E31").Value = "  +5.90%  "
$ws.Range("E32").Value = "  +20.58%  "
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("E34").Value = "  +10.54%  "
$ws.Range("E35").Value = "  +13.13%  "
$ws.Range("E36").Value = "  +4.36%  "
$ws.Range("E37").Value = "  +9.17%  "
$ws.Range("E38").Value = "  +8.24%  "
$ws.Range("E39").Value = "  +5.98%  "
$ws.Range("E40").Value = "  +11.13%  "
$ws.Range("E41").Value = "  +9.43%  "
$ws.Range("E42").Value = "  +4.41%  "
$ws.Range("E43").Value = "  +7.73%  "
$ws.Range("E44").Value = "  +9.70%  "
$ws.Range("E45").Value = "  +2.45%  "
$ws.Range("E46").Value = "  +4.75%  "
$ws.Range("E47").Value = "  +5.87%  "
$ws.Range("E48").Value = "  +6.70%  "
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("E50").Value = "  +6.53%  "
$ws.Range("E51").Value = "  +6.90%  "
